$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Pallindrome" in C7, matching the bold style used by A7
$ws.Range("C7").Value = "Pallindrome"
$ws.Range("A7").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new data rows under the new header
$ws.Range("C8").Value = "Leetcode - 5"
$ws.Range("C9").Value = "Leetcode - 647"

# Match the selection change recorded in the saved workbook
$ws.Range("C7").Select()
